# Update crypto price/volume data per the Dec 14 2023 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value as TEXT (preserves leading/trailing zeros and
# avoids Excel auto-converting numeric-looking strings like "1.00" -> 1).
function Set-TextValue($cell, $text) {
    if ($text.Trim() -match '^[+-]?\d+(\.\d+)?$') {
        $cell.Value = "'" + $text
    } else {
        $cell.Value = $text
    }
}

Set-TextValue $ws.Range("D2") '42.470.85'
Set-TextValue $ws.Range("E2") '  +2.45%  '
Set-TextValue $ws.Range("D3") '2.270.38'
Set-TextValue $ws.Range("E3") '  +3.66%  '
Set-TextValue $ws.Range("E4") '  -0.03%  '
Set-TextValue $ws.Range("D5") '249.93'
Set-TextValue $ws.Range("E5") '  -0.11%  '
Set-TextValue $ws.Range("D6") '0.630'
Set-TextValue $ws.Range("E6") '  +2.52%  '
Set-TextValue $ws.Range("D7") '71.48'
Set-TextValue $ws.Range("E7") '  +5.96%  '
Set-TextValue $ws.Range("E8") '  -0.09%  '
Set-TextValue $ws.Range("D9") '0.641'
Set-TextValue $ws.Range("E9") '  +4.15%  '
Set-TextValue $ws.Range("D10") '38.29'
Set-TextValue $ws.Range("E10") '  -3.48%  '
Set-TextValue $ws.Range("D11") '59.33'
Set-TextValue $ws.Range("E11") '  -0.37%  '
Set-TextValue $ws.Range("D12") '0.0956'
Set-TextValue $ws.Range("E12") '  +1.76%  '
Set-TextValue $ws.Range("D13") '7.28'
Set-TextValue $ws.Range("E13") '  +3.69%  '
Set-TextValue $ws.Range("D14") '0.106'
Set-TextValue $ws.Range("E14") '  +1.76%  '
Set-TextValue $ws.Range("D15") '2.610.63'
Set-TextValue $ws.Range("E15") '  +3.65%  '
Set-TextValue $ws.Range("D16") '14.88'
Set-TextValue $ws.Range("E16") '  +2.66%  '
Set-TextValue $ws.Range("E17") '  +2.42%  '
Set-TextValue $ws.Range("D18") '2.273.55'
Set-TextValue $ws.Range("E18") '  +4.22%  '
Set-TextValue $ws.Range("D19") '42.441.66'
Set-TextValue $ws.Range("E19") '  +2.53%  '
Set-TextValue $ws.Range("D20") '0.0₃0992'
Set-TextValue $ws.Range("E20") '  +4.44%  '
Set-TextValue $ws.Range("D21") '6.26'
Set-TextValue $ws.Range("E21") '  +2.45%  '
Set-TextValue $ws.Range("D22") '72.71'
Set-TextValue $ws.Range("E22") '  +1.08%  '
Set-TextValue $ws.Range("E23") '  +9.57%  '
Set-TextValue $ws.Range("D24") '234.35'
Set-TextValue $ws.Range("E24") '  +1.48%  '
Set-TextValue $ws.Range("E25") '  +1.12%  '
Set-TextValue $ws.Range("D26") '11.52'
Set-TextValue $ws.Range("E26") '  +1.37%  '
Set-TextValue $ws.Range("D27") '1.00'
Set-TextValue $ws.Range("E27") '  -0.02%  '
Set-TextValue $ws.Range("E28") '  +0.74%  '
Set-TextValue $ws.Range("D29") '3.65'
Set-TextValue $ws.Range("E29") '  -0.78%  '
Set-TextValue $ws.Range("E30") '  +1.41%  '
Set-TextValue $ws.Range("D31") '166.62'
Set-TextValue $ws.Range("E31") '  -0.27%  '
Set-TextValue $ws.Range("D32") '20.92'
Set-TextValue $ws.Range("E32") '  +3.11%  '
Set-TextValue $ws.Range("D33") '6.44'
Set-TextValue $ws.Range("E33") '  +11.05%  '
Set-TextValue $ws.Range("D34") '0.127'
Set-TextValue $ws.Range("E34") '  +5.64%  '
Set-TextValue $ws.Range("D35") '31.56'
Set-TextValue $ws.Range("E35") '  +21.28%  '
Set-TextValue $ws.Range("D36") '0.0795'
Set-TextValue $ws.Range("E36") '  +1.67%  '
Set-TextValue $ws.Range("D37") '4.81'
Set-TextValue $ws.Range("E37") '  +13.17%  '
Set-TextValue $ws.Range("E38") '  +3.06%  '
Set-TextValue $ws.Range("E39") '  +3.91%  '
Set-TextValue $ws.Range("E40") '  +0.61%  '
Set-TextValue $ws.Range("D41") '13.60'
Set-TextValue $ws.Range("E41") '  +14.74%  '
Set-TextValue $ws.Range("D42") '2.33'
Set-TextValue $ws.Range("E42") '  +5.15%  '
Set-TextValue $ws.Range("D43") '5.99'
Set-TextValue $ws.Range("E43") '  +6.03%  '
Set-TextValue $ws.Range("D44") '0.209'
Set-TextValue $ws.Range("E44") '  +7.69%  '
Set-TextValue $ws.Range("D45") '9.25'
Set-TextValue $ws.Range("E45") '  +8.16%  '
Set-TextValue $ws.Range("D46") '61.51'
Set-TextValue $ws.Range("E46") '  -0.04%  '
Set-TextValue $ws.Range("D47") '4.88'
Set-TextValue $ws.Range("E47") '  -5.86%  '
Set-TextValue $ws.Range("E48") '  +3.43%  '
Set-TextValue $ws.Range("E49") '  +0.18%  '
Set-TextValue $ws.Range("E50") '  +0.96%  '
Set-TextValue $ws.Range("E51") '  +2.61%  '
